# Bugfix: Checkbox hidden field problem
# Adds three new localization rows (EnterSkill, Abroad, Reference) to both the
# "en" and "de" resource sheets.

$wb   = $excel.ActiveWorkbook
$wsEn = $wb.Worksheets.Item("en")
$wsDe = $wb.Worksheets.Item("de")

# xlPasteFormats - paste only cell formatting (number format/alignment/etc.),
# leaving any existing value untouched. Used so the newly appended rows pick
# up the exact same look as the rows immediately above them without Excel
# fabricating brand-new (unused) cell styles in styles.xml.
$xlPasteFormats = -4122

function Copy-Format($srcWs, $srcAddr, $dstWs, $dstAddr) {
    $srcWs.Range($srcAddr).Copy() | Out-Null
    $dstWs.Range($dstAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---------------------------------------------------------------------------
# Row 204: EnterSkill
#   A204 (both sheets) uses the same "key" style as A203/A202 (centered +
#   wrapped look). Note: on the "de" sheet row 203 happens to use the
#   default format, so A202 (which does carry the style) is used as the
#   source there instead.
#   B204 differs per sheet: "en" keeps the wrapped/value style from B203,
#   "de" uses the default (unstyled) format.
# ---------------------------------------------------------------------------
Copy-Format $wsEn "A203" $wsEn "A204"
Copy-Format $wsDe "A202" $wsDe "A204"
Copy-Format $wsEn "B203" $wsEn "B204"

$wsEn.Range("A204").Value = "EnterSkill"
$wsDe.Range("A204").Value = "EnterSkill"
$wsDe.Range("B204").Value = "Bitte drücken Sie die Enter-Taste um einen neuen Skill einzufügen"
$wsEn.Range("B204").Value = "Please press Enter to insert a new Skill"

# ---------------------------------------------------------------------------
# Row 205: Abroad
#   A205 (both sheets) uses the default (unstyled) format.
#   B205 (both sheets) uses the wrapped/value style from B203.
# ---------------------------------------------------------------------------
Copy-Format $wsEn "B203" $wsEn "B205"
Copy-Format $wsDe "B203" $wsDe "B205"

$wsEn.Range("A205").Value = "Abroad"
$wsDe.Range("A205").Value = "Abroad"
$wsEn.Range("B205").Value = "Abroad stay"
$wsDe.Range("B205").Value = "Auslandaufenthalt"

# ---------------------------------------------------------------------------
# Row 206: Reference
#   A206 (both sheets) uses the same "key" style as A203/A202.
#   B206 (both sheets) uses the wrapped/value style from B203.
# ---------------------------------------------------------------------------
Copy-Format $wsEn "A203" $wsEn "A206"
Copy-Format $wsDe "A202" $wsDe "A206"
Copy-Format $wsEn "B203" $wsEn "B206"
Copy-Format $wsDe "B203" $wsDe "B206"

$wsEn.Range("A206").Value = "Reference"
$wsDe.Range("A206").Value = "Reference"
$wsEn.Range("B206").Value = "References"
$wsDe.Range("B206").Value = "Referenzen"

# ---------------------------------------------------------------------------
# Update the on-screen selection to match the newly added last row, making
# sure the "de" sheet ends up as the active/selected tab (as it was before
# the edit).
# ---------------------------------------------------------------------------
$wsEn.Range("A206:B206").Select() | Out-Null
$wsDe.Range("A206:B206").Select() | Out-Null
